$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (8-10) for the ECs target cluster
$ws.Rows("8:10").Delete()

# Update remaining rows (2-7) with recalculated TPM-based values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Reln"
$ws.Range("C2").Value = "Lrp8"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041827
$ws.Range("H2").Value = 0.125481
$ws.Range("I2").Value = 0.006279874897961605
$ws.Range("J2").Value = 0.006279874897961606
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1124773333333333
$ws.Range("N2").Value = 0.337432
$ws.Range("O2").Value = 0.7871437602495106
$ws.Range("P2").Value = 0.7871437602495107
$ws.Range("Q2").Value = 0.004704589421333333
$ws.Range("R2").Value = 0.042341304792
$ws.Range("S2").Value = 0.004943164341078009
$ws.Range("T2").Value = 0.00494316434107801
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Reln"
$ws.Range("C3").Value = "Lrp8"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.041827
$ws.Range("H3").Value = 0.125481
$ws.Range("I3").Value = 0.006279874897961605
$ws.Range("J3").Value = 0.006279874897961606
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03041566666666666
$ws.Range("N3").Value = 0.09124699999999999
$ws.Range("O3").Value = 0.2128562397504893
$ws.Range("P3").Value = 0.2128562397504893
$ws.Range("Q3").Value = 0.001272196089666667
$ws.Range("R3").Value = 0.011449764807
$ws.Range("S3").Value = 0.001336710556883595
$ws.Range("T3").Value = 0.001336710556883595
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Reln"
$ws.Range("C4").Value = "Lrp8"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.127914
$ws.Range("H4").Value = 3.383742
$ws.Range("I4").Value = 0.1693441751896972
$ws.Range("J4").Value = 0.1693441751896972
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1124773333333333
$ws.Range("N4").Value = 0.337432
$ws.Range("O4").Value = 0.7871437602495106
$ws.Range("P4").Value = 0.7871437602495107
$ws.Range("Q4").Value = 0.1268647589493333
$ws.Range("R4").Value = 1.141782830544
$ws.Range("S4").Value = 0.1332982108351701
$ws.Range("T4").Value = 0.1332982108351702
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Reln"
$ws.Range("C5").Value = "Lrp8"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.127914
$ws.Range("H5").Value = 3.383742
$ws.Range("I5").Value = 0.1693441751896972
$ws.Range("J5").Value = 0.1693441751896972
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03041566666666666
$ws.Range("N5").Value = 0.09124699999999999
$ws.Range("O5").Value = 0.2128562397504893
$ws.Range("P5").Value = 0.2128562397504893
$ws.Range("Q5").Value = 0.03430625625266666
$ws.Range("R5").Value = 0.308756306274
$ws.Range("S5").Value = 0.03604596435452705
$ws.Range("T5").Value = 0.03604596435452705
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Reln"
$ws.Range("C6").Value = "Lrp8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.490742
$ws.Range("H6").Value = 16.472226
$ws.Range("I6").Value = 0.8243759499123412
$ws.Range("J6").Value = 0.8243759499123412
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1124773333333333
$ws.Range("N6").Value = 0.337432
$ws.Range("O6").Value = 0.7871437602495106
$ws.Range("P6").Value = 0.7871437602495107
$ws.Range("Q6").Value = 0.6175840181813332
$ws.Range("R6").Value = 5.558256163631999
$ws.Range("S6").Value = 0.6489023850732625
$ws.Range("T6").Value = 0.6489023850732626
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Reln"
$ws.Range("C7").Value = "Lrp8"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.490742
$ws.Range("H7").Value = 16.472226
$ws.Range("I7").Value = 0.8243759499123412
$ws.Range("J7").Value = 0.8243759499123412
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.03041566666666666
$ws.Range("N7").Value = 0.09124699999999999
$ws.Range("O7").Value = 0.2128562397504893
$ws.Range("P7").Value = 0.2128562397504893
$ws.Range("Q7").Value = 0.1670045784246666
$ws.Range("R7").Value = 1.503041205822
$ws.Range("S7").Value = 0.1754735648390786
$ws.Range("T7").Value = 0.1754735648390787
